# Tab_5b_Wetter.xlsx edit script
# 1) Update row 122 (W_1602x_J -> W_1602x_2030): text, target value, year, and
#    yearly-series markers change; T122 (LabelPositionY) is removed, L122
#    (Zieljahr) is newly added.
# 2) Insert a brand-new row at position 125 (W_1604x_2029), which pushes the
#    former rows 125-128 down to 126-129 (dimension grows from AY128 to AY129).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ThinGridBorder($rng) {
    # Cosmetic helper: restores the thin, light-grey border used throughout
    # this table for any newly-touched cell (matches the surrounding rows).
    for ($i = 7; $i -le 12; $i++) {
        $b = $rng.Borders.Item($i)
        $b.LineStyle = 1
        $b.Weight = 2
        $b.Color = 15132391
    }
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Color = 0
}

# ---------------------------------------------------------------------------
# Row 122: W_1602x_J -> W_1602x_2030
# ---------------------------------------------------------------------------
$ws.Range("A122").Value = "W_1602x_2030"
$ws.Range("G122").Value = "Von 2025 bis 2030 Ausbildung von mindestens 1.000 Personen durch Expertinnen und Experten der Bundeswehr"
$ws.Range("H122").Value = "XXXVon 2025 bis 2030 Ausbildung von mindestens 1.000 Personen durch Expertinnen und Experten der Bundeswehr"
$ws.Range("I122").Value = "K"
$ws.Range("K122").Value = "1000"

# New cell L122 (Zieljahr) = 2030
$ws.Range("L122").Value = 2030
Set-ThinGridBorder $ws.Range("L122")

# T122 (LabelPositionY) is dropped entirely
$ws.Range("T122").ClearContents()

# The yearly trend markers (2015-2021) no longer carry the "S" flag
$ws.Range("AO122").Value = ""
$ws.Range("AP122").Value = ""
$ws.Range("AQ122").Value = ""
$ws.Range("AR122").Value = ""
$ws.Range("AS122").Value = ""
$ws.Range("AT122").Value = ""
$ws.Range("AU122").Value = ""

# ---------------------------------------------------------------------------
# Insert a brand-new row 125 (W_1604x_2029); everything from the old row 125
# onward shifts down by one (125->126, 126->127, 127->128, 128->129).
# ---------------------------------------------------------------------------
$ws.Rows.Item(125).Insert()

$ws.Range("A125").Value = "W_1604x_2029"
$ws.Range("B125").Value = "Z16_B04_P01_IB01_I01"
$ws.Range("C125").Value = "A_BEWERTUNG_PRO"
$ws.Range("D125").Value = $true
$ws.Range("E125").Value = $true
$ws.Range("F125").Value = $true
$ws.Range("G125").Value = "Verbesserung auf mindestens 8 von 11 möglichen Punkten bis 2029"
$ws.Range("H125").Value = "XXXVerbesserung auf mindestens 8 von 11 möglichen Punkten bis 2029"
$ws.Range("I125").Value = "K"
$ws.Range("J125").Value = "steigen"
$ws.Range("K125").Value = "8"
$ws.Range("L125").Value = 2029
$ws.Range("P125").Value = ""
$ws.Range("Q125").Value = ""
$ws.Range("R125").Value = $true
$ws.Range("W125").Value = ""
$ws.Range("X125").Value = ""
$ws.Range("Y125").Value = ""
$ws.Range("Z125").Value = ""
$ws.Range("AA125").Value = ""
$ws.Range("AB125").Value = ""
$ws.Range("AC125").Value = ""
$ws.Range("AD125").Value = ""
$ws.Range("AE125").Value = ""
$ws.Range("AF125").Value = ""
$ws.Range("AG125").Value = ""
$ws.Range("AH125").Value = ""
$ws.Range("AI125").Value = ""
$ws.Range("AJ125").Value = ""
$ws.Range("AK125").Value = ""
$ws.Range("AL125").Value = ""
$ws.Range("AM125").Value = ""
$ws.Range("AN125").Value = ""
$ws.Range("AO125").Value = ""
$ws.Range("AP125").Value = ""
$ws.Range("AQ125").Value = ""
$ws.Range("AR125").Value = ""
$ws.Range("AS125").Value = ""
$ws.Range("AT125").Value = ""
$ws.Range("AU125").Value = ""
$ws.Range("AV125").Value = ""
$ws.Range("AW125").Value = ""
$ws.Range("AX125").Value = ""
$ws.Range("AY125").Value = ""

Set-ThinGridBorder $ws.Range("A125:AY125")

Write-Output "edit complete"
